$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1424943767320741
$ws.Range("D2").Value = 0.3118604158014193
$ws.Range("E2").Value = 0.1840118281489467
$ws.Range("F2").Value = 4.732069065250016
$ws.Range("G2").Value = 0.002378732223644709
$ws.Range("M2").Value = 2.312952077893129
$ws.Range("N2").Value = 1.161793999541331
$ws.Range("B3").Value = 0.1329165875995102
$ws.Range("D3").Value = 0.2791309932840704
$ws.Range("E3").Value = 0.1608958717539366
$ws.Range("F3").Value = 4.290291793991116
$ws.Range("G3").Value = 0.002395785736275016
$ws.Range("M3").Value = 2.037504976932127
$ws.Range("N3").Value = 1.164505009984907
$ws.Range("B4").Value = 0.1271057416350345
$ws.Range("D4").Value = 0.2594259080747179
$ws.Range("E4").Value = 0.1467170950507324
$ws.Range("F4").Value = 4.024650225300491
$ws.Range("G4").Value = 0.002406718192862072
$ws.Range("M4").Value = 1.869911284058247
$ws.Range("N4").Value = 1.166804577164058
$ws.Range("B5").Value = 0.1247554887388134
$ws.Range("D5").Value = 0.2514874601470183
$ws.Range("E5").Value = 0.1409406385789822
$ws.Range("F5").Value = 3.917724988858481
$ws.Range("G5").Value = 0.002411290404933791
$ws.Range("M5").Value = 1.801967288672415
$ws.Range("N5").Value = 1.167902417565884
$ws.Range("B6").Value = 0.1243663046185475
$ws.Range("D6").Value = 0.2501746225975978
$ws.Range("E6").Value = 0.1399814929429581
$ws.Range("F6").Value = 3.900047786284432
$ws.Range("G6").Value = 0.002412056724042516
$ws.Range("M6").Value = 1.79070553704662
$ws.Range("N6").Value = 1.168094453682215
$ws.Range("B7").Value = 0.1270739734135873
$ws.Range("D7").Value = 0.2593184855188326
$ws.Range("E7").Value = 0.1466391883805969
$ws.Range("F7").Value = 4.023202935971142
$ws.Range("G7").Value = 0.002406779379480869
$ws.Range("M7").Value = 1.868993586334454
$ws.Range("N7").Value = 1.166818730722838
$ws.Range("B8").Value = 0.1391775053170647
$ws.Range("D8").Value = 0.3004900153726737
$ws.Range("E8").Value = 0.1760371090530199
$ws.Range("F8").Value = 4.578527167677464
$ws.Range("G8").Value = 0.002384517151032461
$ws.Range("M8").Value = 2.217639334998921
$ws.Range("N8").Value = 1.162597599577566
$ws.Range("B9").Value = 0.1634636790937236
$ws.Range("D9").Value = 0.3846505399327782
$ws.Range("E9").Value = 0.2339055314143508
$ws.Range("F9").Value = 5.716097603498781
$ws.Range("G9").Value = 0.002344471965710315
$ws.Range("M9").Value = 2.915095691529245
$ws.Range("N9").Value = 1.159312783686843
$ws.Range("B10").Value = 0.1816398203814487
$ws.Range("D10").Value = 0.4490420599191509
$ws.Range("E10").Value = 0.2767111162530256
$ws.Range("F10").Value = 6.587449386849585
$ws.Range("G10").Value = 0.002317178501652276
$ws.Range("M10").Value = 3.438292748527005
$ws.Range("N10").Value = 1.159883975067743
$ws.Range("B11").Value = 0.1899805412308098
$ws.Range("D11").Value = 0.4790011395062947
$ws.Range("E11").Value = 0.296284319416813
$ws.Range("F11").Value = 6.992967294441144
$ws.Range("G11").Value = 0.00230520731410648
$ws.Range("M11").Value = 3.679207878160241
$ws.Range("N11").Value = 1.160779577295173
$ws.Range("B12").Value = 0.1931492801008972
$ws.Range("D12").Value = 0.4904513616683914
$ws.Range("E12").Value = 0.3037138658572331
$ws.Range("F12").Value = 7.147960971025498
$ws.Range("G12").Value = 0.002300736687076403
$ws.Range("M12").Value = 3.770902964658546
$ws.Range("N12").Value = 1.161209102450883
$ws.Range("B13").Value = 0.1924663795527408
$ws.Range("D13").Value = 0.4879805089109368
$ws.Range("E13").Value = 0.3021129455562175
$ws.Range("F13").Value = 7.114514612737537
$ws.Range("G13").Value = 0.002301696754342236
$ws.Range("M13").Value = 3.751133289474438
$ws.Range("N13").Value = 1.161112593415183
$ws.Range("B14").Value = 0.1902410296807489
$ws.Range("D14").Value = 0.4799409883821397
$ws.Range("E14").Value = 0.2968951813817569
$ws.Range("F14").Value = 7.005689294240597
$ws.Range("G14").Value = 0.002304838266371068
$ws.Range("M14").Value = 3.686742072815463
$ws.Range("N14").Value = 1.160813109197164
$ws.Range("B15").Value = 0.1888792754276523
$ws.Range("D15").Value = 0.4750305616608728
$ws.Range("E15").Value = 0.2937015380596932
$ws.Range("F15").Value = 6.939220896083327
$ws.Range("G15").Value = 0.00230677064447379
$ws.Range("M15").Value = 3.647362699675881
$ws.Range("N15").Value = 1.160641405737778
$ws.Range("B16").Value = 0.1810961802115827
$ws.Range("D16").Value = 0.447098381875179
$ws.Range("E16").Value = 0.2754342498998881
$ws.Range("F16").Value = 6.561141701914835
$ws.Range("G16").Value = 0.002317969682042819
$ws.Range("M16").Value = 3.422610916122068
$ws.Range("N16").Value = 1.159838149403569
$ws.Range("B17").Value = 0.1763399432845603
$ws.Range("D17").Value = 0.4301402978102828
$ws.Range("E17").Value = 0.2642558114467022
$ws.Range("F17").Value = 6.33162488216044
$ws.Range("G17").Value = 0.00232495289441869
$ws.Range("M17").Value = 3.285511106277994
$ws.Range("N17").Value = 1.159507469544309
$ws.Range("B18").Value = 0.1736110932492494
$ws.Range("D18").Value = 0.4204486461888166
$ws.Range("E18").Value = 0.2578355520029731
$ws.Range("F18").Value = 6.200465240601432
$ws.Range("G18").Value = 0.002329011399124175
$ws.Range("M18").Value = 3.206925848194771
$ws.Range("N18").Value = 1.159377174389661
$ws.Range("B19").Value = 0.172688324805037
$ws.Range("D19").Value = 0.4171776214270722
$ws.Range("E19").Value = 0.2556632651863993
$ws.Range("F19").Value = 6.156199712392606
$ws.Range("G19").Value = 0.002330392783757176
$ws.Range("M19").Value = 3.180363366592672
$ws.Range("N19").Value = 1.159343378920781
$ws.Range("B20").Value = 0.1768455486460709
$ws.Range("D20").Value = 0.4319390116650084
$ws.Range("E20").Value = 0.2654447926668411
$ws.Range("F20").Value = 6.355968281485275
$ws.Range("G20").Value = 0.002324205187653863
$ws.Range("M20").Value = 3.300077241634682
$ws.Range("N20").Value = 1.159536478790898
$ws.Range("B21").Value = 0.1908943904786184
$ws.Range("D21").Value = 0.4822994527801256
$ws.Range("E21").Value = 0.2984272617121206
$ws.Range("F21").Value = 7.037614057614519
$ws.Range("G21").Value = 0.002303913840123058
$ws.Range("M21").Value = 3.70564230273169
$ws.Range("N21").Value = 1.160898630375158
$ws.Range("B22").Value = 0.2001360460211856
$ws.Range("D22").Value = 0.515831862542484
$ws.Range("E22").Value = 0.3200873902742245
$ws.Range("F22").Value = 7.491520990578749
$ws.Range("G22").Value = 0.002291016494449641
$ws.Range("M22").Value = 3.973441295945236
$ws.Range("N22").Value = 1.16231519689147
$ws.Range("B23").Value = 0.1951981517022574
$ws.Range("D23").Value = 0.4978750859357604
$ws.Range("E23").Value = 0.3085163558629773
$ws.Range("F23").Value = 7.248451441533291
$ws.Range("G23").Value = 0.002297867185738176
$ws.Range("M23").Value = 3.830244968102278
$ws.Range("N23").Value = 1.161511333947317
$ws.Range("B24").Value = 0.1766169470917447
$ws.Range("D24").Value = 0.4311256338775138
$ws.Range("E24").Value = 0.264907234637505
$ws.Range("F24").Value = 6.344960170699153
$ws.Range("G24").Value = 0.002324543089701712
$ws.Range("M24").Value = 3.293491162633899
$ws.Range("N24").Value = 1.159523177314554
$ws.Range("B25").Value = 0.1568349367126558
$ws.Range("D25").Value = 0.3614691270802268
$ws.Range("E25").Value = 0.218211561347033
$ws.Range("F25").Value = 5.402556378427278
$ws.Range("G25").Value = 0.002354926125232301
$ws.Range("M25").Value = 2.724703118753609
$ws.Range("N25").Value = 1.159673058208426
